# Generate Report for Archive
#
# 1. Update the "Status" value shown across all three sheets from
#    "Ready for handoff" to "In Translation" (every cell that currently
#    shows "Ready for handoff" gets the new text).
# 2. Narrow the "Status" column on each sheet (Overview: E & F; zh-cn: C;
#    de-de: C) to match the new, shorter status text.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $hit = $used.Find("Ready for handoff")
    if ($hit -ne $null) {
        $first = $hit.Address()
        do {
            $hit.Value = "In Translation"
            $hit = $used.FindNext($hit)
        } while ($hit -ne $null -and $hit.Address() -ne $first)
    }
}

$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5
